$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "wambua"
$ws.Range("B3").Value = "2024-10-22 20:59:02"
